$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Cells.Item(64, 8).Value = 3042.7144
$ws.Cells.Item(64, 9).Value = 2999.75
$ws.Cells.Item(64, 11).Value = 2999.75
$ws.Cells.Item(64, 13).Value = -2751.75
# Row 67
$ws.Cells.Item(67, 8).Value = 3042.7144
$ws.Cells.Item(67, 9).Value = 2999.75
$ws.Cells.Item(67, 11).Value = 2999.75
$ws.Cells.Item(67, 13).Value = -2141.75
# Row 116
$ws.Cells.Item(116, 8).Value = 12678.167
$ws.Cells.Item(116, 9).Value = 27876.25
$ws.Cells.Item(116, 10).Value = 5079.125
$ws.Cells.Item(116, 11).Value = 27876.25
$ws.Cells.Item(116, 12).Value = 5079.125
$ws.Cells.Item(116, 13).Value = -24434.25
$ws.Cells.Item(116, 14).Value = -11963.125
# Row 132
$ws.Cells.Item(132, 8).Value = 1093.3235
$ws.Cells.Item(132, 9).Value = 1005.40625
$ws.Cells.Item(132, 10).Value = 2500
$ws.Cells.Item(132, 11).Value = 3016.21875
$ws.Cells.Item(132, 12).Value = 7500
$ws.Cells.Item(132, 13).Value = -486.21875
$ws.Cells.Item(132, 14).Value = -12560
# Row 137
$ws.Cells.Item(137, 8).Value = 1439.3846
$ws.Cells.Item(137, 9).Value = 1382.1428
$ws.Cells.Item(137, 11).Value = 4146.428400000001
$ws.Cells.Item(137, 13).Value = -1596.428400000001
# Row 138
$ws.Cells.Item(138, 8).Value = 3710.7576
$ws.Cells.Item(138, 9).Value = 3651
$ws.Cells.Item(138, 10).Value = 3774.25
$ws.Cells.Item(138, 11).Value = 10953
$ws.Cells.Item(138, 12).Value = 11322.75
$ws.Cells.Item(138, 13).Value = -5813
$ws.Cells.Item(138, 14).Value = -21602.75

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 4679.3687
$ws.Cells.Item(61, 9).Value = 2886.3572
$ws.Cells.Item(61, 10).Value = 9699.799999999999
$ws.Cells.Item(61, 11).Value = 2886.3572
$ws.Cells.Item(61, 12).Value = 9699.799999999999
$ws.Cells.Item(61, 13).Value = -2674.3572
$ws.Cells.Item(61, 14).Value = -10123.8
# Row 74
$ws.Cells.Item(74, 8).Value = 1273.875
$ws.Cells.Item(74, 9).Value = 855.5405
$ws.Cells.Item(74, 10).Value = 6433.3335
$ws.Cells.Item(74, 11).Value = 855.5405
$ws.Cells.Item(74, 12).Value = 6433.3335
$ws.Cells.Item(74, 13).Value = 18.45950000000005
$ws.Cells.Item(74, 14).Value = -8181.3335
# Row 77
$ws.Cells.Item(77, 8).Value = 1273.875
$ws.Cells.Item(77, 9).Value = 855.5405
$ws.Cells.Item(77, 10).Value = 6433.3335
$ws.Cells.Item(77, 11).Value = 4277.702499999999
$ws.Cells.Item(77, 12).Value = 32166.6675
$ws.Cells.Item(77, 13).Value = 90.29750000000058
$ws.Cells.Item(77, 14).Value = -40902.6675
# Row 97
$ws.Cells.Item(97, 8).Value = 559.7778
$ws.Cells.Item(97, 9).Value = 445.63635
$ws.Cells.Item(97, 10).Value = 739.1429000000001
$ws.Cells.Item(97, 11).Value = 445.63635
$ws.Cells.Item(97, 12).Value = 739.1429000000001
$ws.Cells.Item(97, 13).Value = 50.36365000000001
$ws.Cells.Item(97, 14).Value = -1731.1429
# Row 102
$ws.Cells.Item(102, 8).Value = 1894.7333
$ws.Cells.Item(102, 9).Value = 1647.7693
$ws.Cells.Item(102, 11).Value = 1647.7693
$ws.Cells.Item(102, 13).Value = -25.76929999999993
# Row 122
$ws.Cells.Item(122, 8).Value = 1210.7142
$ws.Cells.Item(122, 9).Value = 1210.7142
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 3632.1426
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -1182.1426
$ws.Cells.Item(122, 14).ClearContents()
# Row 132
$ws.Cells.Item(132, 8).Value = 2067.353
$ws.Cells.Item(132, 9).Value = 1229.125
$ws.Cells.Item(132, 10).Value = 2812.4443
$ws.Cells.Item(132, 11).Value = 3687.375
$ws.Cells.Item(132, 12).Value = 8437.332900000001
$ws.Cells.Item(132, 13).Value = -1157.375
$ws.Cells.Item(132, 14).Value = -13497.3329
# Row 136
$ws.Cells.Item(136, 8).Value = 4679.3687
$ws.Cells.Item(136, 9).Value = 2886.3572
$ws.Cells.Item(136, 10).Value = 9699.799999999999
$ws.Cells.Item(136, 11).Value = 8659.071599999999
$ws.Cells.Item(136, 12).Value = 29099.4
$ws.Cells.Item(136, 13).Value = -6109.071599999999
$ws.Cells.Item(136, 14).Value = -34199.39999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Cells.Item(86, 8).Value = 120193.47
$ws.Cells.Item(86, 9).Value = 2526.3635
$ws.Cells.Item(86, 11).Value = 2526.3635
$ws.Cells.Item(86, 13).Value = -1403.3635
# Row 89
$ws.Cells.Item(89, 8).Value = 120193.47
$ws.Cells.Item(89, 9).Value = 2526.3635
$ws.Cells.Item(89, 11).Value = 12631.8175
$ws.Cells.Item(89, 13).Value = -7015.817499999999
# Row 107
$ws.Cells.Item(107, 8).Value = 3046
$ws.Cells.Item(107, 9).Value = 3046
$ws.Cells.Item(107, 11).Value = 3046
$ws.Cells.Item(107, 13).Value = -1126
# Row 134
$ws.Cells.Item(134, 8).Value = 10993.066
$ws.Cells.Item(134, 9).Value = 10211.72
$ws.Cells.Item(134, 10).Value = 14899.8
$ws.Cells.Item(134, 11).Value = 30635.16
$ws.Cells.Item(134, 12).Value = 44699.39999999999
$ws.Cells.Item(134, 13).Value = -28100.16
$ws.Cells.Item(134, 14).Value = -49769.39999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 10
$ws.Cells.Item(10, 8).Value = 788.44446
$ws.Cells.Item(10, 9).Value = 537
$ws.Cells.Item(10, 10).Value = 2800
$ws.Cells.Item(10, 11).Value = 537
$ws.Cells.Item(10, 12).Value = 2800
$ws.Cells.Item(10, 13).Value = -398
$ws.Cells.Item(10, 14).Value = -3078
# Row 31
$ws.Cells.Item(31, 8).Value = 3572.25
$ws.Cells.Item(31, 9).Value = 2668.6667
$ws.Cells.Item(31, 10).Value = 4734
$ws.Cells.Item(31, 11).Value = 2668.6667
$ws.Cells.Item(31, 12).Value = 4734
$ws.Cells.Item(31, 13).Value = -2373.6667
$ws.Cells.Item(31, 14).Value = -5324
# Row 34
$ws.Cells.Item(34, 8).Value = 3572.25
$ws.Cells.Item(34, 9).Value = 2668.6667
$ws.Cells.Item(34, 10).Value = 4734
$ws.Cells.Item(34, 11).Value = 2668.6667
$ws.Cells.Item(34, 12).Value = 4734
$ws.Cells.Item(34, 13).Value = -2466.6667
$ws.Cells.Item(34, 14).Value = -5138
# Row 132
$ws.Cells.Item(132, 8).Value = 1873.1714
$ws.Cells.Item(132, 9).Value = 1052.619
$ws.Cells.Item(132, 10).Value = 3104
$ws.Cells.Item(132, 11).Value = 3157.857
$ws.Cells.Item(132, 12).Value = 9312
$ws.Cells.Item(132, 13).Value = -627.857
$ws.Cells.Item(132, 14).Value = -14372
# Row 134
$ws.Cells.Item(134, 8).Value = 1158.1666
$ws.Cells.Item(134, 9).Value = 987
$ws.Cells.Item(134, 11).Value = 2961
$ws.Cells.Item(134, 13).Value = -426

$ws = $wb.Worksheets.Item("CUL")
# Row 59
$ws.Cells.Item(59, 8).Value = 5250
$ws.Cells.Item(59, 10).Value = 5250
$ws.Cells.Item(59, 12).Value = 15750
$ws.Cells.Item(59, 14).Value = -16830
# Row 60
$ws.Cells.Item(60, 8).Value = 2965
$ws.Cells.Item(60, 9).Value = 2947.5
$ws.Cells.Item(60, 11).Value = 8842.5
$ws.Cells.Item(60, 13).Value = -8591.5
# Row 61
$ws.Cells.Item(61, 8).Value = 490
$ws.Cells.Item(61, 9).Value = 900
$ws.Cells.Item(61, 10).Value = 285
$ws.Cells.Item(61, 11).Value = 2700
$ws.Cells.Item(61, 12).Value = 855
$ws.Cells.Item(61, 13).Value = -2485
$ws.Cells.Item(61, 14).Value = -1285
# Row 68
$ws.Cells.Item(68, 8).Value = 578.5714
$ws.Cells.Item(68, 9).Value = 578.5714
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 1735.7142
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 13).Value = -924.7142000000001
$ws.Cells.Item(68, 14).ClearContents()
# Row 71
$ws.Cells.Item(71, 8).Value = 578.5714
$ws.Cells.Item(71, 9).Value = 578.5714
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 5207.1426
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 13).Value = -1151.1426
$ws.Cells.Item(71, 14).ClearContents()
# Row 107
$ws.Cells.Item(107, 8).Value = 877.0625
$ws.Cells.Item(107, 10).Value = 877.0625
$ws.Cells.Item(107, 12).Value = 2631.1875
$ws.Cells.Item(107, 14).Value = -6471.1875
# Row 122
$ws.Cells.Item(122, 8).Value = 814.5
$ws.Cells.Item(122, 9).Value = 511.8
$ws.Cells.Item(122, 10).Value = 982.6667
$ws.Cells.Item(122, 11).Value = 4606.2
$ws.Cells.Item(122, 12).Value = 8844.0003
$ws.Cells.Item(122, 13).Value = -2156.2
$ws.Cells.Item(122, 14).Value = -13744.0003
# Row 131
$ws.Cells.Item(131, 8).Value = 748.34
$ws.Cells.Item(131, 10).Value = 775.16486
$ws.Cells.Item(131, 12).Value = 2325.49458
$ws.Cells.Item(131, 14).Value = -12405.49458
# Row 133
$ws.Cells.Item(133, 8).Value = 5000
$ws.Cells.Item(133, 10).Value = 5000
$ws.Cells.Item(133, 12).Value = 15000
$ws.Cells.Item(133, 14).Value = -25120

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Cells.Item(97, 8).Value = 851.73914
$ws.Cells.Item(97, 9).Value = 787.93335
$ws.Cells.Item(97, 10).Value = 971.375
$ws.Cells.Item(97, 11).Value = 787.93335
$ws.Cells.Item(97, 12).Value = 971.375
$ws.Cells.Item(97, 13).Value = -291.93335
$ws.Cells.Item(97, 14).Value = -1963.375
# Row 122
$ws.Cells.Item(122, 8).Value = 2165.4546
$ws.Cells.Item(122, 9).Value = 1762.4
$ws.Cells.Item(122, 10).Value = 2501.3333
$ws.Cells.Item(122, 11).Value = 5287.200000000001
$ws.Cells.Item(122, 12).Value = 7503.999899999999
$ws.Cells.Item(122, 13).Value = -2837.200000000001
$ws.Cells.Item(122, 14).Value = -12403.9999
# Row 132
$ws.Cells.Item(132, 8).Value = 4218.1333
$ws.Cells.Item(132, 9).Value = 3448
$ws.Cells.Item(132, 10).Value = 15000
$ws.Cells.Item(132, 11).Value = 10344
$ws.Cells.Item(132, 12).Value = 45000
$ws.Cells.Item(132, 13).Value = -7814
$ws.Cells.Item(132, 14).Value = -50060

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Cells.Item(46, 8).Value = 1617.7273
$ws.Cells.Item(46, 9).Value = 1033.3334
$ws.Cells.Item(46, 10).Value = 1710
$ws.Cells.Item(46, 11).Value = 1033.3334
$ws.Cells.Item(46, 12).Value = 1710
$ws.Cells.Item(46, 13).Value = -845.3334
$ws.Cells.Item(46, 14).Value = -2086
# Row 55
$ws.Cells.Item(55, 8).Value = 329.80768
$ws.Cells.Item(55, 9).Value = 260.82352
$ws.Cells.Item(55, 10).Value = 460.1111
$ws.Cells.Item(55, 11).Value = 260.82352
$ws.Cells.Item(55, 12).Value = 460.1111
$ws.Cells.Item(55, 13).Value = -87.82351999999997
$ws.Cells.Item(55, 14).Value = -806.1111000000001
# Row 93
$ws.Cells.Item(93, 8).Value = 419.4737
$ws.Cells.Item(93, 9).Value = 368.84616
$ws.Cells.Item(93, 11).Value = 368.84616
$ws.Cells.Item(93, 13).Value = 879.1538399999999
# Row 100
$ws.Cells.Item(100, 8).Value = 1144.75
$ws.Cells.Item(100, 9).Value = 1075.6
$ws.Cells.Item(100, 10).Value = 1260
$ws.Cells.Item(100, 11).Value = 1075.6
$ws.Cells.Item(100, 12).Value = 1260
$ws.Cells.Item(100, 13).Value = -534.5999999999999
$ws.Cells.Item(100, 14).Value = -2342
# Row 132
$ws.Cells.Item(132, 8).Value = 2303.8572
$ws.Cells.Item(132, 9).Value = 2392
$ws.Cells.Item(132, 10).Value = 2237.75
$ws.Cells.Item(132, 11).Value = 7176
$ws.Cells.Item(132, 12).Value = 6713.25
$ws.Cells.Item(132, 13).Value = -4646
$ws.Cells.Item(132, 14).Value = -11773.25
# Row 136
$ws.Cells.Item(136, 8).Value = 3677.5293
$ws.Cells.Item(136, 10).Value = 5166.4443
$ws.Cells.Item(136, 12).Value = 15499.3329
$ws.Cells.Item(136, 14).Value = -20599.3329

$ws = $wb.Worksheets.Item("WVR")
# Row 61
$ws.Cells.Item(61, 8).Value = 29500
$ws.Cells.Item(61, 9).Value = 29000
$ws.Cells.Item(61, 10).Value = 30000
$ws.Cells.Item(61, 11).Value = 29000
$ws.Cells.Item(61, 12).Value = 30000
$ws.Cells.Item(61, 13).Value = -28708
$ws.Cells.Item(61, 14).Value = -30584
# Row 70
$ws.Cells.Item(70, 8).Value = 29400
$ws.Cells.Item(70, 10).Value = 29400
$ws.Cells.Item(70, 12).Value = 29400
$ws.Cells.Item(70, 14).Value = -30030
# Row 73
$ws.Cells.Item(73, 8).Value = 29400
$ws.Cells.Item(73, 10).Value = 29400
$ws.Cells.Item(73, 12).Value = 29400
$ws.Cells.Item(73, 14).Value = -31584
# Row 113
$ws.Cells.Item(113, 8).Value = 577.6429000000001
$ws.Cells.Item(113, 9).Value = 348.91666
$ws.Cells.Item(113, 11).Value = 1046.74998
$ws.Cells.Item(113, 13).Value = 1123.25002
# Row 123
$ws.Cells.Item(123, 8).Value = 45306.062
$ws.Cells.Item(123, 10).Value = 47499.785
$ws.Cells.Item(123, 12).Value = 47499.785
$ws.Cells.Item(123, 14).Value = -57299.785
# Row 132
$ws.Cells.Item(132, 8).Value = 2665.6667
$ws.Cells.Item(132, 9).Value = 2500
$ws.Cells.Item(132, 10).Value = 2698.8
$ws.Cells.Item(132, 11).Value = 7500
$ws.Cells.Item(132, 12).Value = 8096.400000000001
$ws.Cells.Item(132, 13).Value = -4970
$ws.Cells.Item(132, 14).Value = -13156.4
# Row 136
$ws.Cells.Item(136, 8).Value = 2911.32
$ws.Cells.Item(136, 9).Value = 3327.4614
$ws.Cells.Item(136, 11).Value = 9982.3842
$ws.Cells.Item(136, 13).Value = -7432.3842
